$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "286.19"
Set-TextValue $ws.Range("E2") "2.51%"

Set-TextValue $ws.Range("D3") "28.80"
Set-TextValue $ws.Range("E3") "5.03%"

Set-TextValue $ws.Range("D4") "5.034"
Set-TextValue $ws.Range("E4") "4.05%"

Set-TextValue $ws.Range("D5") "0.06706"
Set-TextValue $ws.Range("E5") "5.11%"

Set-TextValue $ws.Range("D6") "7.353"
Set-TextValue $ws.Range("E6") "4.54%"

Set-TextValue $ws.Range("D7") "3.382"
Set-TextValue $ws.Range("E7") "2.27%"

Set-TextValue $ws.Range("D8") "1.372"
Set-TextValue $ws.Range("E8") "6.46%"

Set-TextValue $ws.Range("D9") "0.9426"
Set-TextValue $ws.Range("E9") "5.34%"

Set-TextValue $ws.Range("E10") "1.75%"

Set-TextValue $ws.Range("D11") "0.06703"
Set-TextValue $ws.Range("E11") "16.26%"

Set-TextValue $ws.Range("D12") "0.07549"
Set-TextValue $ws.Range("E12") "1.06%"

Set-TextValue $ws.Range("D13") "0.02974"
Set-TextValue $ws.Range("E13") "2.43%"

Set-TextValue $ws.Range("D14") "0.08998"
Set-TextValue $ws.Range("E14") "0.08%"

Set-TextValue $ws.Range("D15") "0.001602"
Set-TextValue $ws.Range("E15") "2.34%"

Set-TextValue $ws.Range("D16") "0.04483"
Set-TextValue $ws.Range("E16") "1.90%"

Set-TextValue $ws.Range("D17") "0.0006468"
Set-TextValue $ws.Range("E17") "1.07%"

Set-TextValue $ws.Range("D18") "0.006565"
Set-TextValue $ws.Range("E18") "8.38%"

Set-TextValue $ws.Range("D19") "3.495"
Set-TextValue $ws.Range("E19") "0.45%"

Set-TextValue $ws.Range("D20") "2.246"
Set-TextValue $ws.Range("E20") "0.98%"

Set-TextValue $ws.Range("D21") "0.3209"
Set-TextValue $ws.Range("E21") "1.96%"

Set-TextValue $ws.Range("E22") "-3.06%"

Set-TextValue $ws.Range("D23") "4.092"
Set-TextValue $ws.Range("E23") "4.43%"

Set-TextValue $ws.Range("D24") "0.1548"
Set-TextValue $ws.Range("E24") "3.01%"

Set-TextValue $ws.Range("D25") "0.001179"
Set-TextValue $ws.Range("E25") "0.24%"

Set-TextValue $ws.Range("D26") "0.004504"
Set-TextValue $ws.Range("E26") "5.25%"

Set-TextValue $ws.Range("D27") "0.0001246"
Set-TextValue $ws.Range("E27") "5.56%"

Set-TextValue $ws.Range("D28") "0.0001613"
Set-TextValue $ws.Range("E28") "-2.40%"

Set-TextValue $ws.Range("D40") "0.04208"
Set-TextValue $ws.Range("E40") "3.15%"

Set-TextValue $ws.Range("D41") "0.006745"
Set-TextValue $ws.Range("E41") "0.88%"

Set-TextValue $ws.Range("D42") "0.1259"
Set-TextValue $ws.Range("E42") "-10.47%"

Set-TextValue $ws.Range("D43") "0.002013"
Set-TextValue $ws.Range("E43") "-5.05%"

Set-TextValue $ws.Range("E44") "11.53%"

Set-TextValue $ws.Range("D45") "0.00005567"
Set-TextValue $ws.Range("E45") "0.68%"

Set-TextValue $ws.Range("E46") "20.74%"

Set-TextValue $ws.Range("D47") "0.01303"
Set-TextValue $ws.Range("E47") "-29.51%"
